$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seven new form submissions came in since the last export. They belong at
# the top of the log (right below the header row), so push the existing
# rows down by 7 and fill in the new data.
$ws.Rows.Item(2).Resize(7).EntireRow.Insert()

$newRows = @(
    @("parth", "parthpatel082828@gmail.com", 45432.82773920139, "16:55:21", "19:51:46"),
    @("parth", "parthpatel082828@gmail.com", 45429.42889310185, "10:16:54", "10:17:26"),
    @("parth", "parthpatel082828@gmail.com", 45429.425747361114, "10:12:34", "10:12:50"),
    @("parth", "parthpatel082828@gmail.com", 45429.425500335645, "10:12:18", "10:12:28"),
    @("parth", "parthpatel082828@gmail.com", 45429.4253, "14:11:47", "10:12:11"),
    @("mical", "madibic334@nweal.com", 45428.78524579861, "18:49:4", "18:49:49"),
    @("parth", "parthpatel082828@gmail.com", 45427.591582546294, "10:9:1", "14:11:38")
)

# Grab the date-formatted style from an existing date cell so the new cells
# share the same number format (instead of minting brand-new style records).
$ws.Cells.Item(9, 3).Copy()

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$excel.CutCopyMode = $false
